# Weekly update: insert 4 new rows of price data (week of 2021-10-22, serial 44491)
# immediately below the header/first data block (row 586), shifting the existing
# historical rows (587-675) down to (591-679).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 587; this pushes rows 587:675 down to 591:679
$ws.Rows("587:590").Insert()

# Row 587: Tomate, Larga vida, Primera, $/bandeja 18 kilos, Región de Arica y Parinacota
$ws.Cells.Item(587,1).Value2  = 8
$ws.Cells.Item(587,2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(587,3).Value2  = "Coquimbo"
$ws.Cells.Item(587,4).Value2  = 44491
$ws.Cells.Item(587,5).Value2  = 4
$ws.Cells.Item(587,6).Value2  = 100112020
$ws.Cells.Item(587,7).Value2  = "Tomate"
$ws.Cells.Item(587,8).Value2  = "Larga vida"
$ws.Cells.Item(587,9).Value2  = "Primera"
$ws.Cells.Item(587,10).Value2 = 720
$ws.Cells.Item(587,11).Value2 = 13000
$ws.Cells.Item(587,12).Value2 = 14000
$ws.Cells.Item(587,13).Value2 = 13500
$ws.Cells.Item(587,14).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(587,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(587,16).Value2 = 750
$ws.Cells.Item(587,17).Value2 = 18
$ws.Cells.Item(587,18).Value2 = "Hortaliza"

# Row 588: Tomate, Larga vida, Primera, $/caja 10 kilos, Región de Arica y Parinacota
$ws.Cells.Item(588,1).Value2  = 8
$ws.Cells.Item(588,2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(588,3).Value2  = "Coquimbo"
$ws.Cells.Item(588,4).Value2  = 44491
$ws.Cells.Item(588,5).Value2  = 4
$ws.Cells.Item(588,6).Value2  = 100112020
$ws.Cells.Item(588,7).Value2  = "Tomate"
$ws.Cells.Item(588,8).Value2  = "Larga vida"
$ws.Cells.Item(588,9).Value2  = "Primera"
$ws.Cells.Item(588,10).Value2 = 880
$ws.Cells.Item(588,11).Value2 = 6300
$ws.Cells.Item(588,12).Value2 = 6500
$ws.Cells.Item(588,13).Value2 = 6400
$ws.Cells.Item(588,14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(588,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(588,16).Value2 = 640
$ws.Cells.Item(588,17).Value2 = 10
$ws.Cells.Item(588,18).Value2 = "Hortaliza"

# Row 589: Tomate, Larga vida, Segunda, $/bandeja 18 kilos, Región de Arica y Parinacota
$ws.Cells.Item(589,1).Value2  = 8
$ws.Cells.Item(589,2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(589,3).Value2  = "Coquimbo"
$ws.Cells.Item(589,4).Value2  = 44491
$ws.Cells.Item(589,5).Value2  = 4
$ws.Cells.Item(589,6).Value2  = 100112020
$ws.Cells.Item(589,7).Value2  = "Tomate"
$ws.Cells.Item(589,8).Value2  = "Larga vida"
$ws.Cells.Item(589,9).Value2  = "Segunda"
$ws.Cells.Item(589,10).Value2 = 560
$ws.Cells.Item(589,11).Value2 = 11000
$ws.Cells.Item(589,12).Value2 = 12000
$ws.Cells.Item(589,13).Value2 = 11500
$ws.Cells.Item(589,14).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(589,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(589,16).Value2 = 639
$ws.Cells.Item(589,17).Value2 = 18
$ws.Cells.Item(589,18).Value2 = "Hortaliza"

# Row 590: Tomate, Larga vida, Segunda, $/caja 10 kilos, Región de Arica y Parinacota
$ws.Cells.Item(590,1).Value2  = 8
$ws.Cells.Item(590,2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(590,3).Value2  = "Coquimbo"
$ws.Cells.Item(590,4).Value2  = 44491
$ws.Cells.Item(590,5).Value2  = 4
$ws.Cells.Item(590,6).Value2  = 100112020
$ws.Cells.Item(590,7).Value2  = "Tomate"
$ws.Cells.Item(590,8).Value2  = "Larga vida"
$ws.Cells.Item(590,9).Value2  = "Segunda"
$ws.Cells.Item(590,10).Value2 = 560
$ws.Cells.Item(590,11).Value2 = 5000
$ws.Cells.Item(590,12).Value2 = 5500
$ws.Cells.Item(590,13).Value2 = 5250
$ws.Cells.Item(590,14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(590,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(590,16).Value2 = 525
$ws.Cells.Item(590,17).Value2 = 10
$ws.Cells.Item(590,18).Value2 = "Hortaliza"
